$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("pron-poss")
$ws2.Name = "poss"
$ws3 = $wb.Worksheets.Item("prep-coal")
$ws3.Name = "art-coal"
$newSheet = $wb.Worksheets.Add($null, $ws3)
$newSheet.Name = "prep"

$newSheet.Range("A1").Value = "French"
$newSheet.Range("B1").Value = "English"
$newSheet.Range("C1").Value = "French"
$newSheet.Range("D1").Value = "English"
$newSheet.Range("A2").Value = "devant"
$newSheet.Range("B2").Value = "front"
$newSheet.Range("C2").Value = "derrière"
$newSheet.Range("D2").Value = "behind"
$newSheet.Range("A3").Value = "sur"
$newSheet.Range("B3").Value = "on"
$newSheet.Range("C3").Value = "sous"
$newSheet.Range("D3").Value = "under"

foreach ($s in $wb.Worksheets) {
    Write-Output $s.Name
}
